# Update Leve-crafting profit-tracker sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Refreshes the currentAveragePrice* columns (H:L) with new Market Board data and
# recomputes the dependent LeveProfitNQ/HQ columns (M:N) to match.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 4500
$ws.Range("I4").Value = 4000
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 4000
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = -3886
$ws.Range("N4").Value = -5228
$ws.Range("H28").Value = 649.6
$ws.Range("I28").Value = 549.6923
$ws.Range("J28").Value = 1299
$ws.Range("K28").Value = 549.6923
$ws.Range("L28").Value = 1299
$ws.Range("M28").Value = -64.69230000000005
$ws.Range("N28").Value = -2269
$ws.Range("H80").Value = 608.6818
$ws.Range("J80").Value = 886.3333
$ws.Range("L80").Value = 2658.9999
$ws.Range("N80").Value = -4654.9999
$ws.Range("H83").Value = 608.6818
$ws.Range("J83").Value = 886.3333
$ws.Range("L83").Value = 7976.9997
$ws.Range("N83").Value = -17960.9997
$ws.Range("H92").Value = 632.1667
$ws.Range("I92").Value = 265.1111
$ws.Range("K92").Value = 265.1111
$ws.Range("M92").Value = 982.8888999999999
$ws.Range("H107").Value = 552.8
$ws.Range("I107").Value = 678.125
$ws.Range("J107").Value = 51.5
$ws.Range("K107").Value = 678.125
$ws.Range("L107").Value = 51.5
$ws.Range("M107").Value = 1241.875
$ws.Range("N107").Value = -3891.5
$ws.Range("H138").Value = 4166.2
$ws.Range("I138").Value = 3451.25
$ws.Range("K138").Value = 10353.75
$ws.Range("M138").Value = -5213.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1175.8148
$ws.Range("I2").Value = 528.34784
$ws.Range("K2").Value = 528.34784
$ws.Range("M2").Value = -415.34784
$ws.Range("H32").Value = 6395.326
$ws.Range("I32").Value = 4359.775
$ws.Range("J32").Value = 19965.666
$ws.Range("K32").Value = 4359.775
$ws.Range("L32").Value = 19965.666
$ws.Range("M32").Value = -4072.775
$ws.Range("N32").Value = -20539.666
$ws.Range("H74").Value = 954.8125
$ws.Range("I74").Value = 935.13336
$ws.Range("K74").Value = 935.13336
$ws.Range("M74").Value = -61.13336000000004
$ws.Range("H77").Value = 954.8125
$ws.Range("I77").Value = 935.13336
$ws.Range("K77").Value = 4675.6668
$ws.Range("M77").Value = -307.6668
$ws.Range("H110").Value = 3819.818
$ws.Range("I110").Value = 2892
$ws.Range("K110").Value = 2892
$ws.Range("M110").Value = -847
$ws.Range("H116").Value = 1175.8148
$ws.Range("I116").Value = 528.34784
$ws.Range("K116").Value = 528.34784
$ws.Range("M116").Value = 1765.65216

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1175.8148
$ws.Range("I3").Value = 528.34784
$ws.Range("K3").Value = 528.34784
$ws.Range("M3").Value = -414.34784
$ws.Range("H86").Value = 1812.8667
$ws.Range("I86").Value = 1813.1
$ws.Range("J86").Value = 1812.4
$ws.Range("K86").Value = 1813.1
$ws.Range("L86").Value = 1812.4
$ws.Range("M86").Value = -690.0999999999999
$ws.Range("N86").Value = -4058.4
$ws.Range("H89").Value = 1812.8667
$ws.Range("I89").Value = 1813.1
$ws.Range("J89").Value = 1812.4
$ws.Range("K89").Value = 9065.5
$ws.Range("L89").Value = 9062
$ws.Range("M89").Value = -3449.5
$ws.Range("N89").Value = -20294
$ws.Range("H99").Value = 2591.0476
$ws.Range("I99").Value = 2392.875
$ws.Range("J99").Value = 2713
$ws.Range("K99").Value = 2392.875
$ws.Range("L99").Value = 2713
$ws.Range("M99").Value = -894.875
$ws.Range("N99").Value = -5709
$ws.Range("H134").Value = 1958
$ws.Range("J134").Value = 4000
$ws.Range("L134").Value = 12000
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2428.5
$ws.Range("I134").Value = 2232.1875
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 6696.5625
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -4161.5625
$ws.Range("N134").Value = -17067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 1213.8
$ws.Range("I15").Value = 512
$ws.Range("J15").Value = 1681.6666
$ws.Range("K15").Value = 1536
$ws.Range("L15").Value = 5044.9998
$ws.Range("M15").Value = -1396
$ws.Range("N15").Value = -5324.9998
$ws.Range("H60").Value = 269
$ws.Range("I60").Value = 269
$ws.Range("K60").Value = 807
$ws.Range("M60").Value = -556
$ws.Range("H122").Value = 376.08334
$ws.Range("I122").Value = 297.2
$ws.Range("J122").Value = 432.42856
$ws.Range("K122").Value = 2674.8
$ws.Range("L122").Value = 3891.85704
$ws.Range("M122").Value = -224.7999999999997
$ws.Range("N122").Value = -8791.857039999999
$ws.Range("H129").Value = 1205.2858
$ws.Range("I129").Value = 985.75
$ws.Range("J129").Value = 1498
$ws.Range("K129").Value = 2957.25
$ws.Range("L129").Value = 4494
$ws.Range("M129").Value = 2042.75
$ws.Range("N129").Value = -14494
$ws.Range("H130").Value = 5663.3335
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 5663.3335
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 16990.0005
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -27030.0005
$ws.Range("H131").Value = 753
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2883124.5
$ws.Range("I11").Value = 2504249.5
$ws.Range("K11").Value = 2504249.5
$ws.Range("M11").Value = -2504110.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2661
$ws.Range("I7").Value = 2556.125
$ws.Range("K7").Value = 2556.125
$ws.Range("M7").Value = -2444.125
$ws.Range("H13").Value = 5883294
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H20").Value = 1750
$ws.Range("I20").Value = 1750
$ws.Range("K20").Value = 1750
$ws.Range("M20").Value = -1524
$ws.Range("H25").Value = 20000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 20000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 20000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -20460
$ws.Range("H61").Value = 2839.125
$ws.Range("I61").Value = 2695.0667
$ws.Range("K61").Value = 2695.0667
$ws.Range("M61").Value = -2493.0667
$ws.Range("H62").Value = 17294.8
$ws.Range("I62").Value = 12226
$ws.Range("J62").Value = 18562
$ws.Range("K62").Value = 12226
$ws.Range("L62").Value = 18562
$ws.Range("M62").Value = -11602
$ws.Range("N62").Value = -19810
$ws.Range("H65").Value = 17294.8
$ws.Range("I65").Value = 12226
$ws.Range("J65").Value = 18562
$ws.Range("K65").Value = 36678
$ws.Range("L65").Value = 55686
$ws.Range("M65").Value = -33558
$ws.Range("N65").Value = -61926
$ws.Range("H113").Value = 2839.125
$ws.Range("I113").Value = 2695.0667
$ws.Range("K113").Value = 2695.0667
$ws.Range("M113").Value = -525.0666999999999
$ws.Range("H126").Value = 2661
$ws.Range("I126").Value = 2556.125
$ws.Range("K126").Value = 7668.375
$ws.Range("M126").Value = -5198.375
$ws.Range("H132").Value = 4059.1667
$ws.Range("I132").Value = 2872.8667
$ws.Range("K132").Value = 8618.6001
$ws.Range("M132").Value = -6088.6001
$ws.Range("H136").Value = 5183.1875
$ws.Range("I136").Value = 4149.8184
$ws.Range("K136").Value = 12449.4552
$ws.Range("M136").Value = -9899.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 7458
$ws.Range("I24").Value = 3650
$ws.Range("J24").Value = 9996.666999999999
$ws.Range("K24").Value = 3650
$ws.Range("L24").Value = 9996.666999999999
$ws.Range("M24").Value = -3420
$ws.Range("N24").Value = -10456.667
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H30").Value = 9467.5
$ws.Range("J30").Value = 9702.5
$ws.Range("L30").Value = 9702.5
$ws.Range("N30").Value = -9916.5
